$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A24").Copy()
$ws.Range("A25:A26").PasteSpecial(-4122)

$ws.Range("A25").Value = "Alle Templates harmonisieren"
$ws.Range("A26").Value = "Bericht"

$ws.Range("D26").Value = 41518
$ws.Range("D26").NumberFormat = "mm-dd-yy"
$ws.Range("D26").Borders.LineStyle = 1

$ws.Range("D27").Select()
